$d = $word.ActiveDocument

# Locate the paragraph containing the salutation text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Dear Members of the Admissions Committee*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range

    # Replace the text (keep the paragraph mark out of the replaced range).
    $r.Text = "To whom it may concern,"

    # Bold the paragraph mark itself (pPr/rPr) as well as the run.
    $target.Range.Font.Bold = $true

    # Ensure the run(s) font matches the ".AppleSystemUIFont" family used by the
    # rest of the document / the paragraph mark's run properties.
    $target.Range.Font.NameAscii = ".AppleSystemUIFont"
    $target.Range.Font.NameFarEast = ".AppleSystemUIFont"
    $target.Range.Font.NameOther = ".AppleSystemUIFont"
    $target.Range.Font.NameBi = ".AppleSystemUIFont"
}
